$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correção nos dados: os rótulos de cabeçalho "unnamed: 1_level_1" e
# "unnamed: 5_level_1" (artefatos de exportação) são substituídos por
# "total", alinhando com a célula C2 que já continha esse rótulo.
$ws.Range("B2").Value2 = "total"
$ws.Range("F2").Value2 = "total"
